# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between the existing "2021-Q4" sheet
#    and the "总计" (totals) sheet, populated with the quarter's per-fund
#    holdings (same column layout as "2021-Q4").
# 2. Update the "总计" sheet with a new row for 2022-Q1 (kept first) while
#    preserving the existing 2021-Q4 totals row underneath it.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, positioned right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# Re-fetch the "总计" sheet by name now that the new sheet has been
# inserted - worksheet handles track by position, so a reference taken
# before the insert would now (incorrectly) resolve to the new sheet.
$wsTotal = $wb.Worksheets.Item("总计")

# Columns B-G (fund code, name, size, position, ratio, value) hold
# text-formatted numbers in the source data, so force text format
# BEFORE writing them - otherwise Excel auto-converts numeric-looking
# text to real numbers (e.g. dropping the leading zero on fund codes).
$wsQ1.Range("B2:G7").NumberFormat = "@"

# Header row
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

$data = @(
    @(0, "008404", "华泰紫金泰盈混合A", "4.29", "79.73", "4.11", "0.1763", 6),
    @(1, "008405", "华泰紫金泰盈混合C", "3.71", "79.73", "4.11", "0.1525", 6),
    @(2, "005310", "广发电子信息传媒产业精选股票A", "3.99", "90.16", "3.36", "0.1341", 10),
    @(3, "011694", "华泰紫金信息科技主题6个月定期开放混合A", "2.60", "77.49", "4.37", "0.1136", 5),
    @(4, "011695", "华泰紫金信息科技主题6个月定期开放混合C", "0.83", "77.49", "4.37", "0.0363", 5),
    @(5, "010236", "广发电子信息传媒产业精选股票C", "0.81", "90.16", "3.36", "0.0272", 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]
    $wsQ1.Cells.Item($row, 1).Value = $rec[0]
    $wsQ1.Cells.Item($row, 2).Value = $rec[1]
    $wsQ1.Cells.Item($row, 3).Value = $rec[2]
    $wsQ1.Cells.Item($row, 4).Value = $rec[3]
    $wsQ1.Cells.Item($row, 5).Value = $rec[4]
    $wsQ1.Cells.Item($row, 6).Value = $rec[5]
    $wsQ1.Cells.Item($row, 7).Value = $rec[6]
    $wsQ1.Cells.Item($row, 8).Value = $rec[7]
}

# Pull in the header + first-column formatting (bold/centered/bordered
# style) from the "2021-Q4" sheet so the new sheet matches its look.
# Done last, via a formats-only paste, so it doesn't disturb the text
# vs. number typing already locked in above.
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("A2:H2").Copy()
$wsQ1.Range("A2:H7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. "总计" sheet: add the 2022-Q1 total as the new row 2, pushing the
#    existing 2021-Q4 total row down to row 3.
# ---------------------------------------------------------------------

# Copy the index column's style (bold/centered/bordered) down to row 3
# before overwriting row 2/3 contents.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.9

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 6
$wsTotal.Range("D2").Value = 0.64

Write-Output "2022-Q1 sheet added; totals sheet updated"
